$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.73076433333334
$ws.Range("H2").Value = 53.19229300000001
$ws.Range("I2").Value = 0.004631884691211661
$ws.Range("J2").Value = 0.00463188469121166
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.107333666666666
$ws.Range("N2").Value = 21.322001
$ws.Range("O2").Value = 0.7373665550576455
$ws.Range("P2").Value = 0.7373665550576454
$ws.Range("Q2").Value = 126.0184582820326
$ws.Range("R2").Value = 1134.166124538293
$ws.Range("S2").Value = 0.003415396858182988
$ws.Range("T2").Value = 0.003415396858182987
$ws.Range("G3").Value = 17.73076433333334
$ws.Range("H3").Value = 53.19229300000001
$ws.Range("I3").Value = 0.004631884691211661
$ws.Range("J3").Value = 0.00463188469121166
$ws.Range("O3").Value = 0.1688878844614928
$ws.Range("P3").Value = 0.1688878844614928
$ws.Range("Q3").Value = 28.86351527116323
$ws.Range("R3").Value = 259.771637440469
$ws.Range("S3").Value = 0.0007822692065683124
$ws.Range("T3").Value = 0.000782269206568312
$ws.Range("G4").Value = 17.73076433333334
$ws.Range("H4").Value = 53.19229300000001
$ws.Range("I4").Value = 0.004631884691211661
$ws.Range("J4").Value = 0.00463188469121166
$ws.Range("M4").Value = 0.8135026666666666
$ws.Range("N4").Value = 2.440508
$ws.Range("O4").Value = 0.08439869112428164
$ws.Range("P4").Value = 0.08439869112428162
$ws.Range("Q4").Value = 14.42402406720489
$ws.Range("R4").Value = 129.816216604844
$ws.Range("S4").Value = 0.0003909250053768616
$ws.Range("T4").Value = 0.0003909250053768615
$ws.Range("G5").Value = 17.73076433333334
$ws.Range("H5").Value = 53.19229300000001
$ws.Range("I5").Value = 0.004631884691211661
$ws.Range("J5").Value = 0.00463188469121166
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09009266666666667
$ws.Range("N5").Value = 0.270278
$ws.Range("O5").Value = 0.009346869356580103
$ws.Range("P5").Value = 0.009346869356580103
$ws.Range("Q5").Value = 1.597411840828223
$ws.Range("R5").Value = 14.376706567454
$ws.Range("S5").Value = 0.00004329362108349877
$ws.Range("T5").Value = 0.00004329362108349876
$ws.Range("I6").Value = 0.9353873458333681
$ws.Range("J6").Value = 0.935387345833368
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.107333666666666
$ws.Range("N6").Value = 21.322001
$ws.Range("O6").Value = 0.7373665550576455
$ws.Range("P6").Value = 0.7373665550576454
$ws.Range("Q6").Value = 25448.83542591129
$ws.Range("R6").Value = 229039.5188332016
$ws.Range("S6").Value = 0.6897233448416651
$ws.Range("T6").Value = 0.689723344841665
$ws.Range("I7").Value = 0.9353873458333681
$ws.Range("J7").Value = 0.935387345833368
$ws.Range("O7").Value = 0.1688878844614928
$ws.Range("P7").Value = 0.1688878844614928
$ws.Range("R7").Value = 52459.661383467
$ws.Range("S7").Value = 0.1579755899898483
$ws.Range("T7").Value = 0.1579755899898483
$ws.Range("I8").Value = 0.9353873458333681
$ws.Range("J8").Value = 0.935387345833368
$ws.Range("M8").Value = 0.8135026666666666
$ws.Range("N8").Value = 2.440508
$ws.Range("O8").Value = 0.08439869112428164
$ws.Range("P8").Value = 0.08439869112428162
$ws.Range("Q8").Value = 2912.863874625083
$ws.Range("R8").Value = 26215.77487162575
$ws.Range("S8").Value = 0.07894546768255205
$ws.Range("T8").Value = 0.07894546768255202
$ws.Range("I9").Value = 0.9353873458333681
$ws.Range("J9").Value = 0.935387345833368
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09009266666666667
$ws.Range("N9").Value = 0.270278
$ws.Range("O9").Value = 0.009346869356580103
$ws.Range("P9").Value = 0.009346869356580103
$ws.Range("Q9").Value = 322.5898142132369
$ws.Range("R9").Value = 2903.308327919132
$ws.Range("S9").Value = 0.008742943319302704
$ws.Range("T9").Value = 0.008742943319302704
$ws.Range("G10").Value = 227.2177583333333
$ws.Range("H10").Value = 681.653275
$ws.Range("I10").Value = 0.0593570833501536
$ws.Range("J10").Value = 0.05935708335015359
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.107333666666666
$ws.Range("N10").Value = 21.322001
$ws.Range("O10").Value = 0.7373665550576455
$ws.Range("P10").Value = 0.7373665550576454
$ws.Range("Q10").Value = 1614.912423467031
$ws.Range("R10").Value = 14534.21181120327
$ws.Range("S10").Value = 0.04376792806817228
$ws.Range("T10").Value = 0.04376792806817227
$ws.Range("G11").Value = 227.2177583333333
$ws.Range("H11").Value = 681.653275
$ws.Range("I11").Value = 0.0593570833501536
$ws.Range("J11").Value = 0.05935708335015359
$ws.Range("O11").Value = 0.1688878844614928
$ws.Range("P11").Value = 0.1688878844614928
$ws.Range("Q11").Value = 369.8827142608972
$ws.Range("R11").Value = 3328.944428348075
$ws.Range("S11").Value = 0.01002469223481194
$ws.Range("T11").Value = 0.01002469223481194
$ws.Range("G12").Value = 227.2177583333333
$ws.Range("H12").Value = 681.653275
$ws.Range("I12").Value = 0.0593570833501536
$ws.Range("J12").Value = 0.05935708335015359
$ws.Range("M12").Value = 0.8135026666666666
$ws.Range("N12").Value = 2.440508
$ws.Range("O12").Value = 0.08439869112428164
$ws.Range("P12").Value = 0.08439869112428162
$ws.Range("Q12").Value = 184.8422523181889
$ws.Range("R12").Value = 1663.5802708637
$ws.Range("S12").Value = 0.005009660143707854
$ws.Range("T12").Value = 0.005009660143707853
$ws.Range("G13").Value = 227.2177583333333
$ws.Range("H13").Value = 681.653275
$ws.Range("I13").Value = 0.0593570833501536
$ws.Range("J13").Value = 0.05935708335015359
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09009266666666667
$ws.Range("N13").Value = 0.270278
$ws.Range("O13").Value = 0.009346869356580103
$ws.Range("P13").Value = 0.009346869356580103
$ws.Range("Q13").Value = 20.47065376227222
$ws.Range("R13").Value = 184.23588386045
$ws.Range("S13").Value = 0.0005548029034615217
$ws.Range("T13").Value = 0.0005548029034615216
$ws.Range("G14").Value = 2.387458333333333
$ws.Range("H14").Value = 7.162374999999999
$ws.Range("I14").Value = 0.0006236861252666267
$ws.Range("J14").Value = 0.0006236861252666266
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 7.107333666666666
$ws.Range("N14").Value = 21.322001
$ws.Range("O14").Value = 0.7373665550576455
$ws.Range("P14").Value = 0.7373665550576454
$ws.Range("Q14").Value = 16.96846299026389
$ws.Range("R14").Value = 152.716166912375
$ws.Range("S14").Value = 0.0004598852896251037
$ws.Range("T14").Value = 0.0004598852896251036
$ws.Range("G15").Value = 2.387458333333333
$ws.Range("H15").Value = 7.162374999999999
$ws.Range("I15").Value = 0.0006236861252666267
$ws.Range("J15").Value = 0.0006236861252666266
$ws.Range("O15").Value = 0.1688878844614928
$ws.Range("P15").Value = 0.1688878844614928
$ws.Range("Q15").Value = 3.886490100930555
$ws.Range("R15").Value = 34.97841090837499
$ws.Range("S15").Value = 0.0001053330302642662
$ws.Range("T15").Value = 0.0001053330302642662
$ws.Range("G16").Value = 2.387458333333333
$ws.Range("H16").Value = 7.162374999999999
$ws.Range("I16").Value = 0.0006236861252666267
$ws.Range("J16").Value = 0.0006236861252666266
$ws.Range("M16").Value = 0.8135026666666666
$ws.Range("N16").Value = 2.440508
$ws.Range("O16").Value = 0.08439869112428164
$ws.Range("P16").Value = 0.08439869112428162
$ws.Range("Q16").Value = 1.942203720722222
$ws.Range("R16").Value = 17.4798334865
$ws.Range("S16").Value = 0.00005263829264487806
$ws.Range("T16").Value = 0.00005263829264487804
$ws.Range("G17").Value = 2.387458333333333
$ws.Range("H17").Value = 7.162374999999999
$ws.Range("I17").Value = 0.0006236861252666267
$ws.Range("J17").Value = 0.0006236861252666266
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09009266666666667
$ws.Range("N17").Value = 0.270278
$ws.Range("O17").Value = 0.009346869356580103
$ws.Range("P17").Value = 0.009346869356580103
$ws.Range("Q17").Value = 0.2150924878055555
$ws.Range("R17").Value = 1.93583239025
$ws.Range("S17").Value = 0.000005829512732378813
$ws.Range("T17").Value = 0.000005829512732378812
